$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4 timestamp (tiny precision correction from re-scrape)
$ws.Range("A4").Value = 45805.39372784722

# Append new row 5 with the latest scraped price entry
$ws.Range("A5").Value = 45806.39326444409
$ws.Range("B5").Value = "EVOWHEY PROTEIN"
$ws.Range("C5").Value = "2Kg"
$ws.Range("D5").Value = "37,90€"

# Match the date-time style used by the other date cells in column A
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
